$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO")

# Insert two new blank rows at position 5 (pushes the old rows 5-21 down to 7-23)
$ws.Range("A5:A6").EntireRow.Insert()

# Populate the new rows - order matters so new shared-string indices come out
# in the same order as the target workbook (55, 56, 57)
$ws.Range("B5").Value = "Création engagement"
$ws.Range("C6").Value = "Tous les champs doivent être renseignés"
$ws.Range("C5").Value = "Le montant de création doit être > 0 et < au solde restant de la ligne"

# B6 stays empty but still needs to exist as a cell
$ws.Range("B6").Value = ""

# Update the active selection shown in the sheet view
$ws.Activate()
$ws.Range("J3").Select()
